# "Active abilitys implemented without touch buttons"
# Fill in the bullet/ability prompt text that was generated for the sheet,
# bold the small header cell next to the first prompt, scroll the view down
# to where the new content lives, and set the page to print on A4 portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New prompt text (becomes shared strings 0,1,2)
$ws.Range("C37").Value = "bullet0,cosmic,geometric 2D, simple shapes, triangular, symmetric, sharp edges, game asset, vector graphic, solid color background,"
$ws.Range("C39").Value = "shadows, gradient background"
$ws.Range("C41").Value = "bullet,geometric, 2D, simple shapes, triangular, symmetric, sharp edges, game asset, vector graphic, solid color background"

# B37 is an empty "label" cell that just carries the bold style alongside C37
$ws.Range("B37:C37").Font.Bold = $true

# Scroll the window so row 26 is at the top, and leave the selection on C39
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C39").Select()

# Page setup: A4 (paperSize 9), portrait orientation (xlPortrait = 1)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
